$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 403, shifting the existing rows 403:485 down to 406:488.
$ws.Rows("403:405").Insert()

# New weekly price entries for 2023-10-12 (Excel serial 45211), "Cultivar IV Región".
# Row 403: Especial
$ws.Cells.Item(403, 1).Value = 6
$ws.Cells.Item(403, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(403, 3).Value = "Metropolitana"
$ws.Cells.Item(403, 4).Value = 45211
$ws.Cells.Item(403, 5).Value = 13
$ws.Cells.Item(403, 6).Value = "Fruta"
$ws.Cells.Item(403, 7).Value = 100107
$ws.Cells.Item(403, 8).Value = "Otros"
$ws.Cells.Item(403, 9).Value = 100107002
$ws.Cells.Item(403, 10).Value = "Chirimoya"
$ws.Cells.Item(403, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(403, 12).Value = "Especial"
$ws.Cells.Item(403, 13).Value = 278
$ws.Cells.Item(403, 14).Value = 26000
$ws.Cells.Item(403, 15).Value = 26000
$ws.Cells.Item(403, 16).Value = 26000
$ws.Cells.Item(403, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(403, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(403, 19).Value = 2600
$ws.Cells.Item(403, 20).Value = 10

# Row 404: Primera
$ws.Cells.Item(404, 1).Value = 6
$ws.Cells.Item(404, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(404, 3).Value = "Metropolitana"
$ws.Cells.Item(404, 4).Value = 45211
$ws.Cells.Item(404, 5).Value = 13
$ws.Cells.Item(404, 6).Value = "Fruta"
$ws.Cells.Item(404, 7).Value = 100107
$ws.Cells.Item(404, 8).Value = "Otros"
$ws.Cells.Item(404, 9).Value = 100107002
$ws.Cells.Item(404, 10).Value = "Chirimoya"
$ws.Cells.Item(404, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(404, 12).Value = "Primera"
$ws.Cells.Item(404, 13).Value = 350
$ws.Cells.Item(404, 14).Value = 23000
$ws.Cells.Item(404, 15).Value = 23000
$ws.Cells.Item(404, 16).Value = 23000
$ws.Cells.Item(404, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(404, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(404, 19).Value = 2300
$ws.Cells.Item(404, 20).Value = 10

# Row 405: Segunda
$ws.Cells.Item(405, 1).Value = 6
$ws.Cells.Item(405, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(405, 3).Value = "Metropolitana"
$ws.Cells.Item(405, 4).Value = 45211
$ws.Cells.Item(405, 5).Value = 13
$ws.Cells.Item(405, 6).Value = "Fruta"
$ws.Cells.Item(405, 7).Value = 100107
$ws.Cells.Item(405, 8).Value = "Otros"
$ws.Cells.Item(405, 9).Value = 100107002
$ws.Cells.Item(405, 10).Value = "Chirimoya"
$ws.Cells.Item(405, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(405, 12).Value = "Segunda"
$ws.Cells.Item(405, 13).Value = 300
$ws.Cells.Item(405, 14).Value = 20000
$ws.Cells.Item(405, 15).Value = 20000
$ws.Cells.Item(405, 16).Value = 20000
$ws.Cells.Item(405, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(405, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(405, 19).Value = 2000
$ws.Cells.Item(405, 20).Value = 10

# Ensure the date column keeps its date number format for the new rows.
$ws.Range("D403:D405").NumberFormat = $ws.Cells.Item(406, 4).NumberFormat
